$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(138, 8).Value = 3538.5322
$ws.Cells.Item(138, 9).Value = 1781.9131
$ws.Cells.Item(138, 10).Value = 4574.4873
$ws.Cells.Item(138, 11).Value = 5345.7393
$ws.Cells.Item(138, 12).Value = 13723.4619
$ws.Cells.Item(138, 13).Value = -205.7393000000002
$ws.Cells.Item(138, 14).Value = -24003.4619

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 2651.1428
$ws.Cells.Item(2, 9).Value = 2309
$ws.Cells.Item(2, 10).Value = 3506.5
$ws.Cells.Item(2, 11).Value = 2309
$ws.Cells.Item(2, 12).Value = 3506.5
$ws.Cells.Item(2, 13).Value = -2196
$ws.Cells.Item(2, 14).Value = -3732.5
$ws.Cells.Item(32, 8).Value = 12158.322
$ws.Cells.Item(32, 9).Value = 13160.036
$ws.Cells.Item(32, 10).Value = 4287.7144
$ws.Cells.Item(32, 11).Value = 13160.036
$ws.Cells.Item(32, 12).Value = 4287.7144
$ws.Cells.Item(32, 13).Value = -12873.036
$ws.Cells.Item(32, 14).Value = -4861.7144
$ws.Cells.Item(74, 8).Value = 1482.3158
$ws.Cells.Item(74, 9).Value = 1772.7273
$ws.Cells.Item(74, 10).Value = 1083
$ws.Cells.Item(74, 11).Value = 1772.7273
$ws.Cells.Item(74, 12).Value = 1083
$ws.Cells.Item(74, 13).Value = -898.7273
$ws.Cells.Item(74, 14).Value = -2831
$ws.Cells.Item(77, 8).Value = 1482.3158
$ws.Cells.Item(77, 9).Value = 1772.7273
$ws.Cells.Item(77, 10).Value = 1083
$ws.Cells.Item(77, 11).Value = 8863.636500000001
$ws.Cells.Item(77, 12).Value = 5415
$ws.Cells.Item(77, 13).Value = -4495.636500000001
$ws.Cells.Item(77, 14).Value = -14151
$ws.Cells.Item(116, 8).Value = 2651.1428
$ws.Cells.Item(116, 9).Value = 2309
$ws.Cells.Item(116, 10).Value = 3506.5
$ws.Cells.Item(116, 11).Value = 2309
$ws.Cells.Item(116, 12).Value = 3506.5
$ws.Cells.Item(116, 13).Value = -15
$ws.Cells.Item(116, 14).Value = -8094.5
$ws.Cells.Item(122, 8).Value = 2063.3713
$ws.Cells.Item(122, 9).Value = 2110.72
$ws.Cells.Item(122, 11).Value = 6332.16
$ws.Cells.Item(122, 13).Value = -3882.16

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 2651.1428
$ws.Cells.Item(3, 9).Value = 2309
$ws.Cells.Item(3, 10).Value = 3506.5
$ws.Cells.Item(3, 11).Value = 2309
$ws.Cells.Item(3, 12).Value = 3506.5
$ws.Cells.Item(3, 13).Value = -2195
$ws.Cells.Item(3, 14).Value = -3734.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2295.3635
$ws.Cells.Item(31, 9).Value = 1418.5358
$ws.Cells.Item(31, 11).Value = 1418.5358
$ws.Cells.Item(31, 13).Value = -1123.5358
$ws.Cells.Item(34, 8).Value = 2295.3635
$ws.Cells.Item(34, 9).Value = 1418.5358
$ws.Cells.Item(34, 11).Value = 1418.5358
$ws.Cells.Item(34, 13).Value = -1216.5358
$ws.Cells.Item(58, 8).Value = 1982.1724
$ws.Cells.Item(58, 9).Value = 1776.1765
$ws.Cells.Item(58, 10).Value = 2274
$ws.Cells.Item(58, 11).Value = 1776.1765
$ws.Cells.Item(58, 12).Value = 2274
$ws.Cells.Item(58, 13).Value = -1573.1765
$ws.Cells.Item(58, 14).Value = -2680
$ws.Cells.Item(105, 8).Value = 1799.9
$ws.Cells.Item(105, 9).Value = 1766.5
$ws.Cells.Item(105, 10).Value = 1850
$ws.Cells.Item(105, 11).Value = 1766.5
$ws.Cells.Item(105, 12).Value = 1850
$ws.Cells.Item(105, 13).Value = -19.5
$ws.Cells.Item(105, 14).Value = -5344
$ws.Cells.Item(122, 8).Value = 2916.4
$ws.Cells.Item(122, 9).Value = 2933.7144
$ws.Cells.Item(122, 11).Value = 8801.143199999999
$ws.Cells.Item(122, 13).Value = -6351.143199999999
$ws.Cells.Item(136, 8).Value = 1982.1724
$ws.Cells.Item(136, 9).Value = 1776.1765
$ws.Cells.Item(136, 10).Value = 2274
$ws.Cells.Item(136, 11).Value = 5328.529500000001
$ws.Cells.Item(136, 12).Value = 6822
$ws.Cells.Item(136, 13).Value = -2778.529500000001
$ws.Cells.Item(136, 14).Value = -11922

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2146.3
$ws.Cells.Item(5, 9).Value = 2312.6
$ws.Cells.Item(5, 11).Value = 6937.799999999999
$ws.Cells.Item(5, 13).Value = -6825.799999999999
$ws.Cells.Item(57, 8).Value = 1445
$ws.Cells.Item(57, 10).Value = 2000
$ws.Cells.Item(57, 12).Value = 6000
$ws.Cells.Item(57, 14).Value = -7118
$ws.Cells.Item(135, 8).Value = 2146.3
$ws.Cells.Item(135, 9).Value = 2312.6
$ws.Cells.Item(135, 11).Value = 20813.4
$ws.Cells.Item(135, 13).Value = -18278.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3175.5386
$ws.Cells.Item(80, 9).Value = 2822.75
$ws.Cells.Item(80, 10).Value = 3740
$ws.Cells.Item(80, 11).Value = 2822.75
$ws.Cells.Item(80, 12).Value = 3740
$ws.Cells.Item(80, 13).Value = -1824.75
$ws.Cells.Item(80, 14).Value = -5736
$ws.Cells.Item(83, 8).Value = 3175.5386
$ws.Cells.Item(83, 9).Value = 2822.75
$ws.Cells.Item(83, 10).Value = 3740
$ws.Cells.Item(83, 11).Value = 14113.75
$ws.Cells.Item(83, 12).Value = 18700
$ws.Cells.Item(83, 13).Value = -9121.75
$ws.Cells.Item(83, 14).Value = -28684
$ws.Cells.Item(113, 8).Value = 2953.05
$ws.Cells.Item(113, 9).Value = 1544.4286
$ws.Cells.Item(113, 10).Value = 3711.5386
$ws.Cells.Item(113, 11).Value = 1544.4286
$ws.Cells.Item(113, 12).Value = 3711.5386
$ws.Cells.Item(113, 13).Value = 625.5714
$ws.Cells.Item(113, 14).Value = -8051.5386
$ws.Cells.Item(132, 8).Value = 3375.9312
$ws.Cells.Item(132, 9).Value = 1973.091
$ws.Cells.Item(132, 10).Value = 4233.222
$ws.Cells.Item(132, 11).Value = 5919.272999999999
$ws.Cells.Item(132, 12).Value = 12699.666
$ws.Cells.Item(132, 13).Value = -3389.272999999999
$ws.Cells.Item(132, 14).Value = -17759.666

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 8251
$ws.Cells.Item(7, 9).Value = 8251
$ws.Cells.Item(7, 11).Value = 8251
$ws.Cells.Item(7, 13).Value = -8139
$ws.Cells.Item(22, 8).Value = 3924.5
$ws.Cells.Item(22, 9).Value = 4149.3335
$ws.Cells.Item(22, 10).Value = 3250
$ws.Cells.Item(22, 11).Value = 4149.3335
$ws.Cells.Item(22, 12).Value = 3250
$ws.Cells.Item(22, 13).Value = -3854.3335
$ws.Cells.Item(22, 14).Value = -3840
$ws.Cells.Item(27, 8).Value = 3924.5
$ws.Cells.Item(27, 9).Value = 4149.3335
$ws.Cells.Item(27, 10).Value = 3250
$ws.Cells.Item(27, 11).Value = 4149.3335
$ws.Cells.Item(27, 12).Value = 3250
$ws.Cells.Item(27, 13).Value = -4042.3335
$ws.Cells.Item(27, 14).Value = -3464
$ws.Cells.Item(40, 8).Value = 5894.875
$ws.Cells.Item(40, 9).Value = 7893.3335
$ws.Cells.Item(40, 10).Value = 4695.8
$ws.Cells.Item(40, 11).Value = 7893.3335
$ws.Cells.Item(40, 12).Value = 4695.8
$ws.Cells.Item(40, 13).Value = -7757.3335
$ws.Cells.Item(40, 14).Value = -4967.8
$ws.Cells.Item(122, 8).Value = 18185430
$ws.Cells.Item(122, 9).Value = 3887.1667
$ws.Cells.Item(122, 10).Value = 40003280
$ws.Cells.Item(122, 11).Value = 11661.5001
$ws.Cells.Item(122, 12).Value = 120009840
$ws.Cells.Item(122, 13).Value = -9211.500100000001
$ws.Cells.Item(122, 14).Value = -120014740
$ws.Cells.Item(126, 8).Value = 8251
$ws.Cells.Item(126, 9).Value = 8251
$ws.Cells.Item(126, 11).Value = 24753
$ws.Cells.Item(126, 13).Value = -22283

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(107, 8).Value = 769.36365
$ws.Cells.Item(107, 9).Value = 764.7646999999999
$ws.Cells.Item(107, 10).Value = 785
$ws.Cells.Item(107, 11).Value = 2294.2941
$ws.Cells.Item(107, 12).Value = 2355
$ws.Cells.Item(107, 13).Value = -374.2941000000001
$ws.Cells.Item(107, 14).Value = -6195
$ws.Cells.Item(122, 8).Value = 3173.7
$ws.Cells.Item(122, 9).Value = 2246.8
$ws.Cells.Item(122, 10).Value = 3482.6667
$ws.Cells.Item(122, 11).Value = 6740.400000000001
$ws.Cells.Item(122, 12).Value = 10448.0001
$ws.Cells.Item(122, 13).Value = -4290.400000000001
$ws.Cells.Item(122, 14).Value = -15348.0001
$ws.Cells.Item(123, 8).Value = 21534.965
$ws.Cells.Item(123, 10).Value = 21534.965
$ws.Cells.Item(123, 12).Value = 21534.965
$ws.Cells.Item(123, 14).Value = -31334.965
$ws.Cells.Item(126, 8).Value = 10367.777
$ws.Cells.Item(126, 9).Value = 13965.546
$ws.Cells.Item(126, 10).Value = 4714.143
$ws.Cells.Item(126, 11).Value = 41896.638
$ws.Cells.Item(126, 12).Value = 14142.429
$ws.Cells.Item(126, 13).Value = -39426.638
$ws.Cells.Item(126, 14).Value = -19082.429
